{"js": "// Merge the five runs that make up the \"\u4ef7\u683c\uff1a$XXXXXX/X\u4e2a\u6708\uff0c\u4e00\u6b21\u6027\u4ed8\u6b3e\"\n// paragraph into a single run with plain (default) formatting, matching\n// the author's edit which collapsed multiple <w:r> elements (each\n// carrying its own <w:rPr>) into one run with no run-properties at all.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph by its (unique) concatenated text instead of a\n// hard-coded index, so the script is resilient to minor structural\n// differences elsewhere in the document.\nconst targetText = \"\u4ef7\u683c\uff1a$XXXXXX/X\u4e2a\u6708\uff0c\u4e00\u6b21\u6027\u4ed8\u6b3e\";\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === targetText) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find the target paragraph: \" + targetText);\n}\n\n// Remove every run in the paragraph (this also drops their individual\n// <w:rPr> formatting) and then insert the same text back as one single,\n// unformatted run.\ntarget.clear();\nawait context.sync();\n\ntarget.insertText(targetText, Word.InsertLocation.start);\nawait context.sync();\n", "ps1": "# Merge the five runs that make up the \"\u4ef7\u683c\uff1a$XXXXXX/X\u4e2a\u6708\uff0c\u4e00\u6b21\u6027\u4ed8\u6b3e\"\n# paragraph into a single run with plain (default) formatting, matching\n# the author's edit which collapsed multiple runs (each with its own\n# run properties) into one run carrying no run-properties at all.\n\n$d = $word.ActiveDocument\n\n# NOTE: use single-quoted strings throughout so that literal \"$\" signs\n# in the target text (e.g. \"$XXXXXX\") are NOT treated as PowerShell\n# variable interpolation.\n$targetText = '\u4ef7\u683c\uff1a$XXXXXX/X\u4e2a\u6708\uff0c\u4e00\u6b21\u6027\u4ed8\u6b3e'\n\n# Locate the paragraph by its text instead of a hard-coded index, so the\n# script is resilient to minor structural differences elsewhere in the\n# document.\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    $full = $p.Range.Text\n    $trimmed = $full.TrimEnd([char]13, [char]7)\n    if ($trimmed -eq $targetText) {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find the target paragraph: $targetText\"\n}\n\n# Range covering just the paragraph's own text, excluding the trailing\n# paragraph-mark, so assignment below truly replaces (rather than just\n# inserts before) the existing runs.\n$start = $target.Range.Start\n$end = $target.Range.End - 1\n\n# First blank the range out entirely -- this drops every run (and with\n# them, each run's individual <w:rPr> formatting).\n$rng = $d.Range($start, $end)\n$rng.Text = ''\n\n# Then write the merged text back in at the same spot as one brand-new,\n# unformatted run.\n$rng2 = $d.Range($start, $start)\n$rng2.Text = $targetText\n"}
